$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 317645
$ws.Range("D2").Value = 404785147
$ws.Range("C3").Value = 257
$ws.Range("D3").Value = 306979
$ws.Range("C8").Value = 856
$ws.Range("D8").Value = 1259408
$ws.Range("C10").Value = 116345
$ws.Range("D10").Value = 170479391
$ws.Range("C12").Value = 58834
$ws.Range("D12").Value = 84905431
$ws.Range("C16").Value = 3997
$ws.Range("D16").Value = 5672075
$ws.Range("C20").Value = 6548
$ws.Range("D20").Value = 9131339
$ws.Range("C22").Value = 76719
$ws.Range("D22").Value = 95695584
$ws.Range("C27").Value = 286
$ws.Range("D27").Value = 410147
$ws.Range("C28").Value = 32275
$ws.Range("D28").Value = 47247553
$ws.Range("C30").Value = 11386
$ws.Range("D30").Value = 16374789
$ws.Range("C35").Value = 1800
$ws.Range("D35").Value = 2540531
$ws.Range("C36").Value = 96353
$ws.Range("D36").Value = 121305396
$ws.Range("C44").Value = 44133
$ws.Range("D44").Value = 64675603
$ws.Range("C46").Value = 9056
$ws.Range("D46").Value = 12996244
$ws.Range("C51").Value = 2272
$ws.Range("D51").Value = 3169665
$ws.Range("C52").Value = 68448
$ws.Range("D52").Value = 85870468
$ws.Range("C59").Value = 27977
$ws.Range("D59").Value = 41029769
$ws.Range("C62").Value = 10989
$ws.Range("D62").Value = 15890144
$ws.Range("C70").Value = 20316
$ws.Range("D70").Value = 26606591
$ws.Range("C74").Value = 7537
$ws.Range("D74").Value = 11035021
$ws.Range("C76").Value = 5079
$ws.Range("D76").Value = 7374678
$ws.Range("C79").Value = 139548
$ws.Range("D79").Value = 174004762
$ws.Range("C85").Value = 63143
$ws.Range("D85").Value = 92545516
$ws.Range("C88").Value = 29463
$ws.Range("D88").Value = 42619516
$ws.Range("C90").Value = 2723
$ws.Range("D90").Value = 3920852
$ws.Range("C91").Value = 2779
$ws.Range("D91").Value = 3925445
$ws.Range("C92").Value = 32541
$ws.Range("D92").Value = 44078095
$ws.Range("C96").Value = 7863
$ws.Range("D96").Value = 11560856
$ws.Range("C98").Value = 7193
$ws.Range("D98").Value = 10429452
$ws.Range("C102").Value = 8984
$ws.Range("D102").Value = 12473374
$ws.Range("C104").Value = 2266
$ws.Range("D104").Value = 3339302
$ws.Range("C106").Value = 3025
$ws.Range("D106").Value = 4417834
$ws.Range("C109").Value = 178
$ws.Range("D109").Value = 253043
$ws.Range("C110").Value = 139976
$ws.Range("D110").Value = 173057560
$ws.Range("C116").Value = 52359
$ws.Range("D116").Value = 76751867
$ws.Range("C118").Value = 26720
$ws.Range("D118").Value = 38706826
$ws.Range("C119").Value = 1307
$ws.Range("D119").Value = 1788004
$ws.Range("C122").Value = 2220
$ws.Range("D122").Value = 3117138
$ws.Range("C124").Value = 497154
$ws.Range("D124").Value = 655542822
$ws.Range("C126").Value = 210
$ws.Range("D126").Value = 309736
$ws.Range("C131").Value = 205715
$ws.Range("D131").Value = 302397374
$ws.Range("C134").Value = 177987
$ws.Range("D134").Value = 258693003
$ws.Range("C137").Value = 2844
$ws.Range("D137").Value = 3996716
$ws.Range("C139").Value = 6239
$ws.Range("D139").Value = 8812722
$ws.Range("C142").Value = 44091
$ws.Range("D142").Value = 58853334
$ws.Range("C148").Value = 13952
$ws.Range("D148").Value = 20459287
$ws.Range("C149").Value = 3718
$ws.Range("D149").Value = 5360897
$ws.Range("C155").Value = 17394
$ws.Range("D155").Value = 22980924
$ws.Range("C159").Value = 7094
$ws.Range("D159").Value = 10316395
$ws.Range("C161").Value = 4944
$ws.Range("D161").Value = 7114863
$ws.Range("C166").Value = 15728
$ws.Range("D166").Value = 22822012
$ws.Range("C167").Value = 1777
$ws.Range("D167").Value = 2643230
$ws.Range("C170").Value = 53
$ws.Range("D170").Value = 79190
$ws.Range("C171").Value = 86
$ws.Range("D171").Value = 128949
$ws.Range("C172").Value = 86915
$ws.Range("D172").Value = 108707093
$ws.Range("C176").Value = 13
$ws.Range("D176").Value = 19320
$ws.Range("C179").Value = 33634
$ws.Range("D179").Value = 49322182
$ws.Range("C181").Value = 12882
$ws.Range("D181").Value = 18610299
$ws.Range("C183").Value = 1243
$ws.Range("D183").Value = 1740396
$ws.Range("C185").Value = 1629
$ws.Range("D185").Value = 2288635
$ws.Range("C187").Value = 236441
$ws.Range("D187").Value = 293904050
$ws.Range("C193").Value = 869
$ws.Range("D193").Value = 1277845
$ws.Range("C195").Value = 86070
$ws.Range("D195").Value = 126165518
$ws.Range("C198").Value = 32749
$ws.Range("D198").Value = 47131040
$ws.Range("C201").Value = 5096
$ws.Range("D201").Value = 7256836
$ws.Range("C204").Value = 4804
$ws.Range("D204").Value = 6650353
$ws.Range("C207").Value = 261441
$ws.Range("D207").Value = 323540114
$ws.Range("C216").Value = 94555
$ws.Range("D216").Value = 138324414
$ws.Range("C219").Value = 50987
$ws.Range("D219").Value = 73687242
$ws.Range("C222").Value = 4661
$ws.Range("D222").Value = 6544355
$ws.Range("C225").Value = 5650
$ws.Range("D225").Value = 7811237
$ws.Range("C228").Value = 105240
$ws.Range("D228").Value = 131638155
$ws.Range("C229").Value = 75
$ws.Range("D229").Value = 79164
$ws.Range("C235").Value = 49177
$ws.Range("D235").Value = 72039495
$ws.Range("C237").Value = 12255
$ws.Range("D237").Value = 17619577
$ws.Range("C239").Value = 1889
$ws.Range("D239").Value = 2707382
$ws.Range("C241").Value = 2470
$ws.Range("D241").Value = 3453065
$ws.Range("C242").Value = 254907
$ws.Range("D242").Value = 321805729
$ws.Range("C250").Value = 95089
$ws.Range("D250").Value = 139326822
$ws.Range("C253").Value = 64307
$ws.Range("D253").Value = 93183772
$ws.Range("C255").Value = 2402
$ws.Range("D255").Value = 3389701
$ws.Range("C258").Value = 4534
$ws.Range("D258").Value = 6365450
